# Updates cryptos list values (price/volume) and reorders a few coin rows,
# matching the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''91.784.73'
$ws.Range('E2').Value = '''  +0.56%  '
# Row 3
$ws.Range('D3').Value = '''3.087.86'
$ws.Range('E3').Value = '''  -2.13%  '
# Row 4
$ws.Range('E4').Value = '''  -0.03%  '
# Row 5
$ws.Range('D5').Value = '''232.72'
$ws.Range('E5').Value = '''  -4.23%  '
# Row 6
$ws.Range('D6').Value = '''609.31'
$ws.Range('E6').Value = '''  -1.40%  '
# Row 7
$ws.Range('D7').Value = '''1.09'
$ws.Range('E7').Value = '''  -3.45%  '
# Row 8
$ws.Range('D8').Value = '''0.382'
$ws.Range('E8').Value = '''  +2.01%  '
# Row 9
$ws.Range('E9').Value = '''  -0.11%  '
# Row 10
$ws.Range('D10').Value = '''3.083.55'
$ws.Range('E10').Value = '''  -2.20%  '
# Row 11
$ws.Range('D11').Value = '''0.764'
$ws.Range('E11').Value = '''  +3.30%  '
# Row 12
$ws.Range('D12').Value = '''0.197'
$ws.Range('E12').Value = '''  -3.44%  '
# Row 13
$ws.Range('D13').Value = '''0.0000242'
$ws.Range('E13').Value = '''  -2.22%  '
# Row 14
$ws.Range('D14').Value = '''91.931.05'
$ws.Range('E14').Value = '''  +1.18%  '
# Row 15
$ws.Range('B15').Value = '''Avalanche'
$ws.Range('C15').Value = '''https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '''33.36'
$ws.Range('E15').Value = '''  -5.09%  '
# Row 16
$ws.Range('B16').Value = '''Toncoin'
$ws.Range('C16').Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D16').Value = '''5.37'
$ws.Range('E16').Value = '''  -5.18%  '
# Row 17
$ws.Range('D17').Value = '''3.661.91'
$ws.Range('E17').Value = '''  -2.12%  '
# Row 18
$ws.Range('D18').Value = '''3.079.76'
$ws.Range('E18').Value = '''  -1.95%  '
# Row 19
$ws.Range('E19').Value = '''  +1.82%  '
# Row 20
$ws.Range('D20').Value = '''14.33'
$ws.Range('E20').Value = '''  -5.21%  '
# Row 21
$ws.Range('D21').Value = '''5.76'
$ws.Range('E21').Value = '''  -3.45%  '
# Row 22
$ws.Range('D22').Value = '''434.83'
$ws.Range('E22').Value = '''  -5.20%  '
# Row 23
$ws.Range('D23').Value = '''9.03'
$ws.Range('E23').Value = '''  -1.20%  '
# Row 24
$ws.Range('D24').Value = '''0.0000193'
$ws.Range('E24').Value = '''  -6.60%  '
# Row 25
$ws.Range('D25').Value = '''5.55'
$ws.Range('E25').Value = '''  -6.62%  '
# Row 26
$ws.Range('D26').Value = '''85.44'
$ws.Range('E26').Value = '''  -4.04%  '
# Row 27
$ws.Range('D27').Value = '''11.28'
$ws.Range('E27').Value = '''  -5.35%  '
# Row 28
$ws.Range('D28').Value = '''3.260.40'
$ws.Range('E28').Value = '''  -2.01%  '
# Row 29
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '''  -0.11%  '
# Row 30
$ws.Range('B30').Value = '''Cronos'
$ws.Range('C30').Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D30').Value = '''0.177'
$ws.Range('E30').Value = '''  +4.73%  '
# Row 31
$ws.Range('B31').Value = '''Hedera'
$ws.Range('C31').Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.127'
$ws.Range('E31').Value = '''  -14.58%  '
# Row 32
$ws.Range('D32').Value = '''0.235'
$ws.Range('E32').Value = '''  +0.42%  '
# Row 33
$ws.Range('B33').Value = '''InternetComputer(DFINITY)'
$ws.Range('C33').Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '''9.08'
$ws.Range('E33').Value = '''  -3.32%  '
# Row 34
$ws.Range('B34').Value = '''Binance-PegBSC-USD'
$ws.Range('C34').Value = '''https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = '''0.999'
$ws.Range('E34').Value = '''  +42.23%  '
# Row 35
$ws.Range('D35').Value = '''7.70'
$ws.Range('E35').Value = '''  +1.72%  '
# Row 36
$ws.Range('D36').Value = '''0.157'
$ws.Range('E36').Value = '''  -10.09%  '
# Row 37
$ws.Range('D37').Value = '''25.35'
$ws.Range('E37').Value = '''  -4.48%  '
# Row 38
$ws.Range('D38').Value = '''3.87'
$ws.Range('E38').Value = '''  +0.19%  '
# Row 39
$ws.Range('E39').Value = '''  -3.44%  '
# Row 40
$ws.Range('E40').Value = '''  +7.60%  '
# Row 41
$ws.Range('D41').Value = '''0.433'
$ws.Range('E41').Value = '''  -2.97%  '
# Row 42
$ws.Range('B42').Value = '''Bittensor'
$ws.Range('C42').Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '''464.42'
$ws.Range('E42').Value = '''  -6.83%  '
# Row 43
$ws.Range('B43').Value = '''Fetch.AI'
$ws.Range('C43').Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').Value = '''1.26'
$ws.Range('E43').Value = '''  -5.29%  '
# Row 44
$ws.Range('E44').Value = '''  -7.00%  '
# Row 45
$ws.Range('E45').Value = '''  +0.00%  '
# Row 46
$ws.Range('D46').Value = '''161.09'
$ws.Range('E46').Value = '''  +3.67%  '
# Row 47
$ws.Range('D47').Value = '''0.675'
$ws.Range('E47').Value = '''  -5.55%  '
# Row 48
$ws.Range('D48').Value = '''1.81'
$ws.Range('E48').Value = '''  -6.02%  '
# Row 49
$ws.Range('D49').Value = '''1.32'
$ws.Range('E49').Value = '''  -3.62%  '
# Row 50
$ws.Range('B50').Value = '''OKB'
$ws.Range('C50').Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '''43.81'
$ws.Range('E50').Value = '''  -0.52%  '
# Row 51
$ws.Range('B51').Value = '''FirstDigitalUSD'
$ws.Range('C51').Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').Value = '''0.997'
$ws.Range('E51').Value = '''  -0.08%  '
